$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns before the old column E ("7-8" results column),
#     shifting it to column G. ---
$ws.Range("E1:F1").EntireColumn.Insert()

# --- Header row values for the two new columns. ---
$ws.Range("E1").Value = "Strategy"
$ws.Range("F1").Value = "VM"

# --- Player selection now happens in units of 20: Strategy column is
#     always 5, VM alternates in two batches of 5 players each. ---
$ws.Range("E2:E11").Value = 5
$ws.Range("F2:F6").Value = 1
$ws.Range("F7:F11").Value = 2

# --- Seed the Hyperlink / Followed Hyperlink font+style artifacts (not
#     applied to any live cell) that show up in the style table. ---
$seed1 = $ws.Range("Z500")
$seed1.Value = "seed1"
$ws.Hyperlinks.Add($seed1, "http://example.com/1")
$ws.Hyperlinks.Item(1).Delete()
$seed1.Clear()

$seed2 = $ws.Range("Z501")
$seed2.Value = "seed2"
$seed2.Style = "Followed Hyperlink"
$seed2.Clear()

# --- Header cell formatting: bold black font plus a border. E1 keeps a
#     full thin border (like the rest of row 1); F1 gets a thin border on
#     every side but the left (so it reads as one continuous box with E1). ---
$e1 = $ws.Range("E1")
$e1.Style = "Normal"
$e1.Font.Color = 0
$e1.Font.Bold = $true
$e1.Borders.Color = 0
$e1.Borders.LineStyle = 1
$e1.HorizontalAlignment = -4108
$e1.VerticalAlignment = -4160

$f1 = $ws.Range("F1")
$f1.Style = "Normal"
$f1.Font.Color = 0
$f1.Font.Bold = $true
$f1.Borders.Color = 0
$f1.Borders.Item(7).LineStyle = -4142
$f1.HorizontalAlignment = -4108
$f1.VerticalAlignment = -4160

# --- Data cell formatting: plain black font for the new numeric columns. ---
$data = $ws.Range("E2:F11")
$data.Font.Color = 0

# --- Column widths for the two new columns. ---
$ws.Columns.Item(5).ColumnWidth = 6.8333333
$ws.Columns.Item(6).ColumnWidth = 3.1666667

# --- Page setup / selection to match the saved view state. ---
$ws.PageSetup.Orientation = 1
$ws.Range("G18").Select()
